# Update cryptocurrency price (Price) and volume change (Volume(1h)) columns
# for rows 2-51 on the active worksheet, reflecting refreshed market data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.246.85'
$ws.Range("E2").Value = '  -0.09%  '

$ws.Range("D3").Value = '1.602.05'
$ws.Range("E3").Value = '  -0.39%  '

$ws.Range("E4").Value = '  +0.08%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '212.50'
$ws.Range("E5").Value = '  -0.09%  '

$ws.Range("E6").Value = '  +0.06%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.485'
$ws.Range("E7").Value = '  +0.31%  '

$ws.Range("E8").Value = '  -0.20%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0614'
$ws.Range("E9").Value = '  -0.74%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '18.19'
$ws.Range("E10").Value = '  -0.60%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0814'
$ws.Range("E11").Value = '  +1.05%  '

$ws.Range("D12").Value = '1.822.97'
$ws.Range("E12").Value = '  -0.23%  '

$ws.Range("D13").Value = '1.600.73'
$ws.Range("E13").Value = '  -0.12%  '

$ws.Range("E14").Value = '  +0.66%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.514'
$ws.Range("E15").Value = '  +0.53%  '

$ws.Range("D16").Value = '26.216.92'
$ws.Range("E16").Value = '  -0.04%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '61.37'
$ws.Range("E17").Value = '  +1.01%  '

$ws.Range("D18").Value = '0.0₃0730'
$ws.Range("E18").Value = '  -0.05%  '

$ws.Range("E19").Value = '  +0.02%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '203.05'
$ws.Range("E20").Value = '  +1.62%  '

$ws.Range("E21").Value = '  +0.59%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.27'
$ws.Range("E22").Value = '  -1.60%  '

$ws.Range("E23").Value = '  -0.43%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.91'
$ws.Range("E24").Value = '  +9.12%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '144.12'
$ws.Range("E25").Value = '  +1.28%  '

$ws.Range("E26").Value = '  +0.18%  '

$ws.Range("E27").Value = '  -8.02%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.20'
$ws.Range("E28").Value = '  +0.13%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.55'
$ws.Range("E29").Value = '  +0.72%  '

$ws.Range("E30").Value = '  +3.12%  '

$ws.Range("E31").Value = '  -0.66%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.17'
$ws.Range("E32").Value = '  +1.24%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.93'
$ws.Range("E33").Value = '  -2.91%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.43'
$ws.Range("E34").Value = '  +3.01%  '

$ws.Range("E35").Value = '  -1.52%  '

$ws.Range("D36").Value = '1.159.57'
$ws.Range("E36").Value = '  +4.44%  '

$ws.Range("E37").Value = '  +8.24%  '

$ws.Range("E38").Value = '  +0.10%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.793'
$ws.Range("E39").Value = '  +0.50%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.33'
$ws.Range("E40").Value = '  -1.06%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.498'
$ws.Range("E41").Value = '  -0.90%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.780'
$ws.Range("E42").Value = '  -0.31%  '

$ws.Range("E43").Value = '  +1.87%  '

$ws.Range("D44").Value = '1.737.20'
$ws.Range("E44").Value = '  -0.13%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '91.87'
$ws.Range("E45").Value = '  -1.02%  '

$ws.Range("E46").Value = '  -2.93%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '54.07'
$ws.Range("E47").Value = '  +0.71%  '

$ws.Range("E48").Value = '  -0.63%  '

$ws.Range("E49").Value = '  -0.56%  '

$ws.Range("D50").Value = '0.0₇0958'
$ws.Range("E50").Value = '  -11.26%  '

$ws.Range("E51").Value = '  +0.02%  '
